$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2022" column (S), mirroring the formatting of the existing
# "2021" column (R) immediately to its left (same header style / same
# data-row style as the rest of the year series).
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 76.1

# Match the author's final selection left in the saved workbook.
$ws.Range("P8").Select() | Out-Null
